$d = $word.ActiveDocument
$find = $d.Content.Find

for ($i = 1; $i -le 4; $i++) {
    $find.ClearFormatting()
    $old = "<id>p084v_a$i</id>"
    $new = "<id>p084v_$i</id>"
    $find.Execute($old, $false, $false, $false, $false, $false, $true, 1, $false, $new, 2) | Out-Null
}
